$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L mirroring column K (the 2020 data column),
# copying formatting + values from K3:K4 to L3:L4.
$ws.Range("K3:K4").Copy($ws.Range("L3:L4"))

# Update the selected/active cell as recorded in the sheet view
$ws.Range("L10").Select()
